{"js": "// Update the date heading and the 25 two-digit \u00f7 one-digit division\n// problems in the practice-sheet table. Cells are addressed positionally\n// (table row/column) rather than by searching for their old text, because\n// several of the new values coincide with OTHER cells' old values\n// elsewhere in the table (e.g. row0/col2 becomes \"65\u00f79=7, 2\", which is the\n// ORIGINAL text of row3/col2) \u2014 a content search after earlier edits would\n// otherwise hit the wrong, already-updated cell.\n\nconst body = context.document.body;\n\n// --- 1. Date heading -------------------------------------------------\nconst headingPara = body.paragraphs.getFirst();\nheadingPara.load(\"text\");\nawait context.sync();\n\nif (headingPara.text.trim() === \"2025-03-10 Monday\") {\n  headingPara.getRange().insertText(\"2025-03-11 Tuesday\", \"Replace\");\n} else {\n  // Fallback: search-based replace if the structure differs from expectation.\n  const hits = body.search(\"2025-03-10 Monday\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(\"2025-03-11 Tuesday\", \"Replace\");\n  }\n}\n\n// --- 2. Division problems in the table -------------------------------\nconst table = body.tables.getFirst();\n\n// (row, col, newText) triples, in document order.\nconst cellUpdates = [\n  [0, 0, \"64\u00f75=12, 4\"],\n  [0, 1, \"57\u00f72=28, 1\"],\n  [0, 2, \"65\u00f79=7, 2\"],\n  [0, 3, \"62\u00f79=6, 8\"],\n  [0, 4, \"92\u00f74=23, 0\"],\n  [4, 0, \"60\u00f74=15, 0\"],\n  [4, 1, \"78\u00f77=11, 1\"],\n  [4, 2, \"79\u00f77=11, 2\"],\n  [4, 3, \"49\u00f73=16, 1\"],\n  [4, 4, \"78\u00f79=8, 6\"],\n  [8, 0, \"93\u00f78=11, 5\"],\n  [8, 1, \"76\u00f75=15, 1\"],\n  [8, 2, \"45\u00f75=9, 0\"],\n  [8, 3, \"36\u00f78=4, 4\"],\n  [8, 4, \"25\u00f78=3, 1\"],\n  [12, 0, \"66\u00f74=16, 2\"],\n  [12, 1, \"42\u00f74=10, 2\"],\n  [12, 2, \"11\u00f75=2, 1\"],\n  [12, 3, \"67\u00f75=13, 2\"],\n  [12, 4, \"48\u00f78=6, 0\"],\n  [16, 0, \"43\u00f76=7, 1\"],\n  [16, 1, \"68\u00f72=34, 0\"],\n  [16, 2, \"69\u00f72=34, 1\"],\n  [16, 3, \"62\u00f72=31, 0\"],\n  [16, 4, \"29\u00f78=3, 5\"],\n];\n\nfor (const [row, col, newText] of cellUpdates) {\n  const cell = table.getCell(row, col);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 two-digit / one-digit division\n# problems in the practice-sheet table. Cells are addressed positionally\n# (1-based Table.Cell(row, col)) instead of by searching for their old\n# text, because several NEW values coincide with OTHER cells' OLD values\n# elsewhere in the table (e.g. row1/col3 becomes \"65\u00f79=7, 2\", which is the\n# ORIGINAL text of row4/col3) \u2014 a Find/Replace pass after earlier edits\n# would otherwise hit the wrong, already-updated cell.\n\n$d = $word.ActiveDocument\n\n# --- 1. Date heading ---------------------------------------------------\n$heading = $d.Paragraphs(1).Range\nif ($heading.Text.TrimEnd(\"`r\") -eq \"2025-03-10 Monday\") {\n    $heading.Text = \"2025-03-11 Tuesday\"\n} else {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = \"2025-03-10 Monday\"\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = \"2025-03-11 Tuesday\"\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\n# --- 2. Division problems in the table ----------------------------------\n$t = $d.Tables(1)\n\n# (row, col, newText) triples, 1-based, in document order.\n$cellUpdates = @(\n    @(1, 1, \"64\u00f75=12, 4\"),\n    @(1, 2, \"57\u00f72=28, 1\"),\n    @(1, 3, \"65\u00f79=7, 2\"),\n    @(1, 4, \"62\u00f79=6, 8\"),\n    @(1, 5, \"92\u00f74=23, 0\"),\n    @(5, 1, \"60\u00f74=15, 0\"),\n    @(5, 2, \"78\u00f77=11, 1\"),\n    @(5, 3, \"79\u00f77=11, 2\"),\n    @(5, 4, \"49\u00f73=16, 1\"),\n    @(5, 5, \"78\u00f79=8, 6\"),\n    @(9, 1, \"93\u00f78=11, 5\"),\n    @(9, 2, \"76\u00f75=15, 1\"),\n    @(9, 3, \"45\u00f75=9, 0\"),\n    @(9, 4, \"36\u00f78=4, 4\"),\n    @(9, 5, \"25\u00f78=3, 1\"),\n    @(13, 1, \"66\u00f74=16, 2\"),\n    @(13, 2, \"42\u00f74=10, 2\"),\n    @(13, 3, \"11\u00f75=2, 1\"),\n    @(13, 4, \"67\u00f75=13, 2\"),\n    @(13, 5, \"48\u00f78=6, 0\"),\n    @(17, 1, \"43\u00f76=7, 1\"),\n    @(17, 2, \"68\u00f72=34, 0\"),\n    @(17, 3, \"69\u00f72=34, 1\"),\n    @(17, 4, \"62\u00f72=31, 0\"),\n    @(17, 5, \"29\u00f78=3, 5\")\n)\n\nforeach ($update in $cellUpdates) {\n    $row = $update[0]\n    $col = $update[1]\n    $newText = $update[2]\n    $cell = $t.Cell($row, $col)\n    $cell.Range.Text = $newText\n}\n"}
